# Improve formatting of seating plan excel
# Replace legacy programme abbreviations with the full "B.Tech ..." names
# used throughout the Programme column (and wherever the same short code
# appears as a course code), matching whole-cell contents only so that
# strings like "ITVC" are not clobbered by the "IT" replacement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPart = 1
$xlWhole = 2

# Order matters: replace the longer/more specific token ("ITVC") before the
# shorter one ("IT") so a whole-cell match behaves as expected either way.
$ws.Cells.Replace("Comps", "B.Tech COMP", $xlPart, $xlWhole) | Out-Null
$ws.Cells.Replace("ITVC", "B.Tech ITVC", $xlPart, $xlWhole) | Out-Null
$ws.Cells.Replace("IT", "B.Tech IT", $xlPart, $xlWhole) | Out-Null
